# Update LR-pair NATMI edge-weight metrics (Rarres2 -> Gpr1) with refreshed TPM-derived
# expression values. Target-cluster labels (column D) are unaffected by this refresh;
# only the per-row quantitative columns (E:T) change, driven by updated ligand/receptor
# expression inputs (columns E:J and K:P) which in turn drive the derived edge metrics
# (Q:T = G*M, H*N, I*O, J*P respectively).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data for A2:T21 (20 data rows x 20 columns), rebuilt in full per row so the
# arithmetic relationships between the Ligand/Receptor inputs and the Edge outputs
# stay internally consistent.
$rows = @(
    ,@("ECs", "Rarres2", "Gpr1", "ECs", 2.0, 0.6666666666666666, 0.5310463333333334, 1.593139, 0.01208304439515528, 0.01241653163470256, 2.0, 0.6666666666666666, 0.9956586666666668, 2.986976, 0.2287590061179932, 0.2450602367332084, 0.5287408841848891, 4.758667957664001, 0.002764105226715312, 0.003042798181805581)
    ,@("ECs", "Rarres2", "Gpr1", "FAPs", 2.0, 0.6666666666666666, 0.5310463333333334, 1.593139, 0.01208304439515528, 0.01241653163470256, 3.0, 1.0, 2.226167333333333, 6.678502, 0.51147631580469, 0.5479238136306437, 1.182197999753111, 10.639781997778, 0.006180191030938534, 0.006803313365351758)
    ,@("ECs", "Rarres2", "Gpr1", "Inflammatory-Mac", 2.0, 0.6666666666666666, 0.5310463333333334, 1.593139, 0.01208304439515528, 0.01241653163470256, 2.0, 0.6666666666666666, 0.1761506666666667, 0.528452, 0.04047175280319151, 0.04335574581855945, 0.09354416564755558, 0.8418974908280001, 0.0004890219858707133, 0.0005383279895022668)
    ,@("ECs", "Rarres2", "Gpr1", "MuSCs", 2.0, 0.6666666666666666, 0.5310463333333334, 1.593139, 0.01208304439515528, 0.01241653163470256, 2.0, 1.0, 0.8685624999999999, 1.737125, 0.1995578413599836, 0.1425188095703395, 0.4612469308958334, 2.767481585375, 0.002411266256554038, 0.001769589307570271)
    ,@("ECs", "Rarres2", "Gpr1", "Resolving-Mac", 2.0, 0.6666666666666666, 0.5310463333333334, 1.593139, 0.01208304439515528, 0.01241653163470256, 1.0, 0.3333333333333333, 0.08589566666666666, 0.257687, 0.0197350839141417, 0.02114139424724881, 0.04561457883255556, 0.410531209493, 0.0002384598950766891, 0.0002625027904726836)
    ,@("FAPs", "Rarres2", "Gpr1", "ECs", 3.0, 1.0, 39.49382433333334, 118.481473, 0.8986139302737502, 0.9234153188332328, 2.0, 0.6666666666666666, 0.9956586666666668, 2.986976, 0.2287590061179932, 0.2450602367332084, 39.32236847729423, 353.9013162956481, 0.2055660295732068, 0.2262923766363431)
    ,@("FAPs", "Rarres2", "Gpr1", "FAPs", 3.0, 1.0, 39.49382433333334, 118.481473, 0.8986139302737502, 0.9234153188332328, 3.0, 1.0, 2.226167333333333, 6.678502, 0.51147631580469, 0.5479238136306437, 87.91986159927178, 791.2787543934461, 0.4596197423871903, 0.5059612430600616)
    ,@("FAPs", "Rarres2", "Gpr1", "Inflammatory-Mac", 3.0, 1.0, 39.49382433333334, 118.481473, 0.8986139302737502, 0.9234153188332328, 2.0, 0.6666666666666666, 0.1761506666666667, 0.528452, 0.04047175280319151, 0.04335574581855945, 6.95686348553289, 62.61177136979601, 0.03636848085154359, 0.04003535984829767)
    ,@("FAPs", "Rarres2", "Gpr1", "MuSCs", 3.0, 1.0, 39.49382433333334, 118.481473, 0.8986139302737502, 0.9234153188332328, 2.0, 1.0, 0.8685624999999999, 1.737125, 0.1995578413599836, 0.1425188095703395, 34.30285479752083, 205.817128785125, 0.1793254561414404, 0.1316040519791279)
    ,@("FAPs", "Rarres2", "Gpr1", "Resolving-Mac", 3.0, 1.0, 39.49382433333334, 118.481473, 0.8986139302737502, 0.9234153188332328, 1.0, 0.3333333333333333, 0.08589566666666666, 0.257687, 0.0197350839141417, 0.02114139424724881, 3.392348370327889, 30.531135332951, 0.01773422132036914, 0.01952228730940234)
    ,@("Inflammatory-Mac", "Rarres2", "Gpr1", "ECs", 2.0, 1.0, 3.541247, 7.082494000000001, 0.08057497440313205, 0.0551992078554295, 2.0, 0.6666666666666666, 0.9956586666666668, 2.986976, 0.2287590061179932, 0.2450602367332084, 3.525873266357334, 21.155239598144, 0.01843225106244323, 0.01352713094453713)
    ,@("Inflammatory-Mac", "Rarres2", "Gpr1", "FAPs", 2.0, 1.0, 3.541247, 7.082494000000001, 0.08057497440313205, 0.0551992078554295, 3.0, 1.0, 2.226167333333333, 6.678502, 0.51147631580469, 0.5479238136306437, 7.883408390664667, 47.300450343988, 0.04121219105377118, 0.03024496047753751)
    ,@("Inflammatory-Mac", "Rarres2", "Gpr1", "Inflammatory-Mac", 2.0, 1.0, 3.541247, 7.082494000000001, 0.08057497440313205, 0.0551992078554295, 2.0, 0.6666666666666666, 0.1761506666666667, 0.528452, 0.04047175280319151, 0.04335574581855945, 0.6237930198813334, 3.742758119288001, 0.003261010446167043, 0.002393202825165832)
    ,@("Inflammatory-Mac", "Rarres2", "Gpr1", "MuSCs", 2.0, 1.0, 3.541247, 7.082494000000001, 0.08057497440313205, 0.0551992078554295, 2.0, 1.0, 0.8685624999999999, 1.737125, 0.1995578413599836, 0.1425188095703395, 3.0757943474375, 12.30317738975, 0.01607936795952497, 0.007866925392781547)
    ,@("Inflammatory-Mac", "Rarres2", "Gpr1", "Resolving-Mac", 2.0, 1.0, 3.541247, 7.082494000000001, 0.08057497440313205, 0.0551992078554295, 1.0, 0.3333333333333333, 0.08589566666666666, 0.257687, 0.0197350839141417, 0.02114139424724881, 0.3041777718963333, 1.825066631378, 0.001590153881225631, 0.001166988215407469)
    ,@("MuSCs", "Rarres2", "Gpr1", "ECs", 1.0, 0.3333333333333333, 0.3835953333333333, 1.150786, 0.00872805092796245, 0.008968941676635136, 2.0, 0.6666666666666666, 0.9956586666666668, 2.986976, 0.2287590061179932, 0.2450602367332084, 0.3819300181262223, 3.437370163136001, 0.001996620255627919, 0.002197930970522545)
    ,@("MuSCs", "Rarres2", "Gpr1", "FAPs", 1.0, 0.3333333333333333, 0.3835953333333333, 1.150786, 0.00872805092796245, 0.008968941676635136, 3.0, 1.0, 2.226167333333333, 6.678502, 0.51147631580469, 0.5479238136306437, 0.8539474002857778, 7.685526602572001, 0.004464191332789939, 0.004914296727692743)
    ,@("MuSCs", "Rarres2", "Gpr1", "Inflammatory-Mac", 1.0, 0.3333333333333333, 0.3835953333333333, 1.150786, 0.00872805092796245, 0.008968941676635136, 2.0, 0.6666666666666666, 0.1761506666666667, 0.528452, 0.04047175280319151, 0.04335574581855945, 0.06757057369688889, 0.6081351632720001, 0.0003532395196101625, 0.0003888551555936774)
    ,@("MuSCs", "Rarres2", "Gpr1", "MuSCs", 1.0, 0.3333333333333333, 0.3835953333333333, 1.150786, 0.00872805092796245, 0.008968941676635136, 2.0, 1.0, 0.8685624999999999, 1.737125, 0.1995578413599836, 0.1425188095703395, 0.3331765217083333, 1.99905913025, 0.001741751002464188, 0.001278242890859845)
    ,@("MuSCs", "Rarres2", "Gpr1", "Resolving-Mac", 1.0, 0.3333333333333333, 0.3835953333333333, 1.150786, 0.00872805092796245, 0.008968941676635136, 1.0, 0.3333333333333333, 0.08589566666666666, 0.257687, 0.0197350839141417, 0.02114139424724881, 0.03294917688688889, 0.296542591982, 0.0001896159319663242, 0.00018961593196632415)
)

$data = New-Object 'object[,]' $rows.Count, 20
for ($i = 0; $i -lt $rows.Count; $i++) {
    for ($j = 0; $j -lt 20; $j++) {
        $data[$i, $j] = $rows[$i][$j]
    }
}

$ws.Range("A2:T21").Value = $data
